# Updated cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-formatted cells so numeric-looking strings (e.g. "255.49")
# are stored as text, matching the original inlineStr cell type.
$cellUpdates = @{
    "D2" = "36.722.68"
    "E2" = "  +0.35%  "
    "D3" = "2.126.95"
    "E3" = "  +10.82%  "
    "E4" = "  +0.04%  "
    "D5" = "255.49"
    "E5" = "  +2.70%  "
    "E6" = "  -3.91%  "
    "E7" = "  -0.02%  "
    "D8" = "47.26"
    "E8" = "  +6.90%  "
    "D9" = "59.67"
    "E9" = "  +1.95%  "
    "E10" = "  +1.82%  "
    "D11" = "0.0741"
    "E11" = "  -3.19%  "
    "D13" = "2.436.76"
    "E13" = "  +10.70%  "
    "D14" = "14.37"
    "E14" = "  -1.47%  "
    "D15" = "0.842"
    "E15" = "  +5.61%  "
    "D16" = "2.130.23"
    "E16" = "  +10.93%  "
    "D17" = "5.12"
    "E17" = "  -0.10%  "
    "D18" = "36.738.89"
    "E18" = "  +0.27%  "
    "D19" = "73.65"
    "E19" = "  -0.83%  "
    "D20" = "0.0₃0838"
    "E20" = "  -2.57%  "
    "D21" = "13.35"
    "E21" = "  +0.86%  "
    "D22" = "241.64"
    "E22" = "  -3.72%  "
    "D23" = "5.20"
    "E23" = "  +0.07%  "
    "E24" = "  +0.16%  "
    "D25" = "2.49"
    "E25" = "  -7.76%  "
    "D26" = "171.89"
    "E26" = "  +2.54%  "
    "D27" = "21.60"
    "E27" = "  +15.13%  "
    "D28" = "9.21"
    "E28" = "  +4.44%  "
    "D29" = "2.04"
    "E29" = "  -7.01%  "
    "D30" = "28.27"
    "E30" = "  +59.85%  "
    "D31" = "0.123"
    "E31" = "  -4.25%  "
    "E32" = "  -0.67%  "
    "D33" = "0.0948"
    "E33" = "  +11.33%  "
    "D34" = "0.0601"
    "E34" = "  -1.33%  "
    "D35" = "2.36"
    "E35" = "  +17.45%  "
    "B36" = "ImmutableX"
    "C36" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "D36" = "0.948"
    "E36" = "  +7.79%  "
    "B37" = "WEMIXToken"
    "C37" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D37" = "1.89"
    "E37" = "  -4.72%  "
    "E38" = "  +0.05%  "
    "D39" = "4.17"
    "E39" = "  -3.97%  "
    "E40" = "  -9.39%  "
    "E41" = "  +8.10%  "
    "E42" = "  -1.40%  "
    "D43" = "99.11"
    "E43" = "  -7.41%  "
    "D44" = "2.81"
    "E44" = "  +12.97%  "
    "D45" = "16.13"
    "E45" = "  -6.00%  "
    "D46" = "1.360.73"
    "E46" = "  +1.60%  "
    "D47" = "7.19"
    "E47" = "  +11.89%  "
    "D48" = "0.0841"
    "E48" = "  +3.57%  "
    "D49" = "2.327.90"
    "E49" = "  +10.45%  "
    "B50" = "MXToken"
    "C50" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D50" = "2.84"
    "E50" = "  +1.35%  "
    "B51" = "RenderToken"
    "C51" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D51" = "2.29"
    "E51" = "  -3.12%  "
}

foreach ($ref in $cellUpdates.Keys) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $cellUpdates[$ref]
}